$d = $word.ActiveDocument

# --- 1. Remove the "Meta description" paragraph (it sat right after the
#        Heading1 title at the top of the document). ---
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Meta description*") {
        $p.Range.Delete()
    }
}

# --- 2. Insert a new bold paragraph "Play Book of Santa for Free - Enjoy
#        this Holiday Slot Game" right before the final (italic, image
#        prompt) paragraph. ---
$count = $d.Paragraphs.Count
$lastP = $d.Paragraphs.Item($count)
$insertPoint = $d.Range($lastP.Range.Start, $lastP.Range.Start)
$titleText = "Play Book of Santa for Free - Enjoy this Holiday Slot Game"
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>' + $titleText + '</w:t></w:r></w:p>'
$insertPoint.InsertXML($xml)

# Split what is now one merged paragraph into two: the new bold title
# paragraph, and the original (now pushed-down) last paragraph.
$splitPos = $lastP.Range.Start + $titleText.Length
$splitPoint = $d.Range($splitPos, $splitPos)
$splitPoint.InsertParagraphAfter()

# --- 3. Replace the text of the final paragraph (formerly the "Create a
#        feature image..." AI-art prompt) with the meta-description text,
#        keeping its italic run formatting intact. ---
$oldText = "Create a feature image for Book of Santa that captures the essence of the game's theme and unique character. The image should be in a cartoon style and feature a happy-looking Maya warrior wearing glasses, as the protagonist of the game. The Maya warrior should be holding a large book in his hands, which should resemble Santa's book, with a few gifts spilling out of it. The background of the image should be a cozy fireplace scene, with the grid of the slot game superimposed on it. The image should use bright and cheerful colors and convey the festive mood of the holiday season."
$newText = "Read a review of Book of Santa slot game, a captivating game featuring Santa Claus and his book. Learn how you can play it for free."
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

Write-Output "done"
